$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the existing row 216. This shifts the
# current rows 216-337 down to 218-339 and brings along their formatting
# (the date-formatted style on column D gets inherited automatically).
$ws.Rows("216:217").Insert()

# --- Fill in the new row 216 (a new weekly price record) -----------------
$ws.Range("A216").Value = 4
$ws.Range("B216").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C216").Value = "Los Lagos"
$ws.Range("D216").Value2 = 44767
$ws.Range("E216").Value = 10
$ws.Range("F216").Value = "Fruta"
$ws.Range("G216").Value = 100102
$ws.Range("H216").Value = "Cítricos"
$ws.Range("I216").Value = 100102006
$ws.Range("J216").Value = "Pomelo"
$ws.Range("K216").Value = "Start Ruby"
$ws.Range("L216").Value = "Primera"
$ws.Range("M216").Value = 15
$ws.Range("N216").Value = 15000
$ws.Range("O216").Value = 15000
$ws.Range("P216").Value = 15000
$ws.Range("Q216").Value = "$/caja 14 kilos empedrada"
$ws.Range("R216").Value = "Región de O'Higgins"
$ws.Range("S216").Value = 1071
$ws.Range("T216").Value = 14

# --- Fill in the new row 217 (a new weekly price record) -----------------
$ws.Range("A217").Value = 4
$ws.Range("B217").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C217").Value = "Los Lagos"
$ws.Range("D217").Value2 = 44767
$ws.Range("E217").Value = 10
$ws.Range("F217").Value = "Fruta"
$ws.Range("G217").Value = 100102
$ws.Range("H217").Value = "Cítricos"
$ws.Range("I217").Value = 100102006
$ws.Range("J217").Value = "Pomelo"
$ws.Range("K217").Value = "Start Ruby"
$ws.Range("L217").Value = "Segunda"
$ws.Range("M217").Value = 15
$ws.Range("N217").Value = 12000
$ws.Range("O217").Value = 12000
$ws.Range("P217").Value = 12000
$ws.Range("Q217").Value = "$/caja 14 kilos empedrada"
$ws.Range("R217").Value = "Región de O'Higgins"
$ws.Range("S217").Value = 857
$ws.Range("T217").Value = 14
